$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.487.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.72%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.992.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.62%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("E8").Value = "  +1.37%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.985.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.88%  "

# Row 10
$ws.Range("E10").Value = "  +4.48%  "

# Row 11
$ws.Range("E11").Value = "  +11.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.89%  "

# Row 13
$ws.Range("E13").Value = "  +4.21%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "

# Row 15
$ws.Range("E15").Value = "  -0.38%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.490.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.21%  "

# Row 17
$ws.Range("E17").Value = "  +4.18%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.990.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.476.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "435.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.38%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.43%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "

# Row 26
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("B27").Value = "FirstDigitalUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.04%  "

# Row 29
$ws.Range("E29").Value = "  +2.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.34%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.25%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.40%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0779"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.73%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.992"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.53%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.11%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.02%  "

# Row 40
$ws.Range("E40").Value = "  +6.86%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "401.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.86%  "

# Row 42
$ws.Range("E42").Value = "  +2.81%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.758.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.47%  "

# Row 44
$ws.Range("E44").Value = "  -1.96%  "

# Row 45
$ws.Range("E45").Value = "  +5.71%  "

# Row 46
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "

# Row 49
$ws.Range("E49").Value = "  +1.60%  "

# Row 50
$ws.Range("E50").Value = "  +2.34%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.29%  "
